# Apply edit: add 2022-Q1 fund-holdings sheet, repurposing the old
# "总计" sheet (sheetId 6) for it, and append a brand-new "总计"
# totals sheet (sheetId 7) with the 2022-Q1 row prepended.

function Set-TextCell {
    # Force a cell to hold TEXT even when the content parses as a number
    # (fund codes with leading zeros, decimal strings that must stay literal).
    # NumberFormat="@" + Value prevents Excel auto-coercing to a number; the
    # operation leaves a stray text-format style behind, so reset with
    # Style = "Normal" to keep the cell style-free (matches sibling sheets).
    param($ws, $row, $col, $text)
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Repurpose the existing "总计" sheet (sheetId 6) into "2022-Q1":
#    rename it and replace its contents with the fund-holdings table.
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("总计")
$q1.Name = "2022-Q1"
$q1.Cells.Clear()

# Header row, styled like the other quarterly sheets (style copied from the
# still-intact header of the "2021-Q4" sheet, which carries cellXfs index 2).
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Count; $i++) {
    $q1.Cells.Item(1, $i + 2).Value = $headers[$i]
}
$q4 = $wb.Worksheets.Item("2021-Q4")
$q4.Range("B1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)  # xlPasteFormats

# 2022-Q1 row 2 (index 0)
$q1.Cells.Item(2, 1).Value = 0
Set-TextCell $q1 2 2 '161032'
Set-TextCell $q1 2 3 '富国中证煤炭指数'
Set-TextCell $q1 2 4 '21.64'
Set-TextCell $q1 2 5 '93.99'
Set-TextCell $q1 2 6 '3.22'
Set-TextCell $q1 2 7 '0.6968'
$q1.Cells.Item(2, 8).Value = 10

# 2022-Q1 row 3 (index 1)
$q1.Cells.Item(3, 1).Value = 1
Set-TextCell $q1 3 2 '010338'
Set-TextCell $q1 3 3 '国投瑞银远见成长混合A'
Set-TextCell $q1 3 4 '13.17'
Set-TextCell $q1 3 5 '77.42'
Set-TextCell $q1 3 6 '3.72'
Set-TextCell $q1 3 7 '0.4899'
$q1.Cells.Item(3, 8).Value = 2

# 2022-Q1 row 4 (index 2)
$q1.Cells.Item(4, 1).Value = 2
Set-TextCell $q1 4 2 '168204'
Set-TextCell $q1 4 3 '中融中证煤炭指数'
Set-TextCell $q1 4 4 '8.68'
Set-TextCell $q1 4 5 '92.62'
Set-TextCell $q1 4 6 '3.16'
Set-TextCell $q1 4 7 '0.2743'
$q1.Cells.Item(4, 8).Value = 10

# 2022-Q1 row 5 (index 3)
$q1.Cells.Item(5, 1).Value = 3
Set-TextCell $q1 5 2 '519198'
Set-TextCell $q1 5 3 '万家颐和灵活配置混合'
Set-TextCell $q1 5 4 '1.78'
Set-TextCell $q1 5 5 '91.03'
Set-TextCell $q1 5 6 '8.15'
Set-TextCell $q1 5 7 '0.1451'
$q1.Cells.Item(5, 8).Value = 1

# 2022-Q1 row 6 (index 4)
$q1.Cells.Item(6, 1).Value = 4
Set-TextCell $q1 6 2 '550001'
Set-TextCell $q1 6 3 '信诚四季红混合'
Set-TextCell $q1 6 4 '5.00'
Set-TextCell $q1 6 5 '72.84'
Set-TextCell $q1 6 6 '2.58'
Set-TextCell $q1 6 7 '0.1290'
$q1.Cells.Item(6, 8).Value = 9

# 2022-Q1 row 7 (index 5)
$q1.Cells.Item(7, 1).Value = 5
Set-TextCell $q1 7 2 '121010'
Set-TextCell $q1 7 3 '国投瑞银瑞源灵活配置混合'
Set-TextCell $q1 7 4 '2.99'
Set-TextCell $q1 7 5 '67.79'
Set-TextCell $q1 7 6 '3.12'
Set-TextCell $q1 7 7 '0.0933'
$q1.Cells.Item(7, 8).Value = 1

# 2022-Q1 row 8 (index 6)
$q1.Cells.Item(8, 1).Value = 6
Set-TextCell $q1 8 2 '010339'
Set-TextCell $q1 8 3 '国投瑞银远见成长混合C'
Set-TextCell $q1 8 4 '2.02'
Set-TextCell $q1 8 5 '77.42'
Set-TextCell $q1 8 6 '3.72'
Set-TextCell $q1 8 7 '0.0751'
$q1.Cells.Item(8, 8).Value = 2

# 2022-Q1 row 9 (index 7)
$q1.Cells.Item(9, 1).Value = 7
Set-TextCell $q1 9 2 '161222'
Set-TextCell $q1 9 3 '国投瑞银瑞利灵活配置混合（LOF）'
Set-TextCell $q1 9 4 '2.21'
Set-TextCell $q1 9 5 '63.82'
Set-TextCell $q1 9 6 '3.26'
Set-TextCell $q1 9 7 '0.0720'
$q1.Cells.Item(9, 8).Value = 1

# 2022-Q1 row 10 (index 8)
$q1.Cells.Item(10, 1).Value = 8
Set-TextCell $q1 10 2 '000165'
Set-TextCell $q1 10 3 '国投瑞银策略精选混合'
Set-TextCell $q1 10 4 '3.38'
Set-TextCell $q1 10 5 '70.38'
Set-TextCell $q1 10 6 '2.03'
Set-TextCell $q1 10 7 '0.0686'
$q1.Cells.Item(10, 8).Value = 10

# 2022-Q1 row 11 (index 9)
$q1.Cells.Item(11, 1).Value = 9
Set-TextCell $q1 11 2 '510170'
Set-TextCell $q1 11 3 '国联安上证大宗商品股票ETF'
Set-TextCell $q1 11 4 '2.22'
Set-TextCell $q1 11 5 '98.14'
Set-TextCell $q1 11 6 '2.48'
Set-TextCell $q1 11 7 '0.0551'
$q1.Cells.Item(11, 8).Value = 7

# 2022-Q1 row 12 (index 10)
$q1.Cells.Item(12, 1).Value = 10
Set-TextCell $q1 12 2 '001601'
Set-TextCell $q1 12 3 '鑫元鑫新收益灵活配置混合A'
Set-TextCell $q1 12 4 '0.72'
Set-TextCell $q1 12 5 '86.90'
Set-TextCell $q1 12 6 '4.39'
Set-TextCell $q1 12 7 '0.0316'
$q1.Cells.Item(12, 8).Value = 7

# 2022-Q1 row 13 (index 11)
$q1.Cells.Item(13, 1).Value = 11
Set-TextCell $q1 13 2 '001266'
Set-TextCell $q1 13 3 '国投瑞银招财灵活配置混合'
Set-TextCell $q1 13 4 '0.76'
Set-TextCell $q1 13 5 '67.47'
Set-TextCell $q1 13 6 '3.83'
Set-TextCell $q1 13 7 '0.0291'
$q1.Cells.Item(13, 8).Value = 1

# 2022-Q1 row 14 (index 12)
$q1.Cells.Item(14, 1).Value = 12
Set-TextCell $q1 14 2 '161715'
Set-TextCell $q1 14 3 '招商中证大宗商品股票指数（LOF）'
Set-TextCell $q1 14 4 '2.24'
Set-TextCell $q1 14 5 '95.07'
Set-TextCell $q1 14 6 '1.24'
Set-TextCell $q1 14 7 '0.0278'
$q1.Cells.Item(14, 8).Value = 9

# 2022-Q1 row 15 (index 13)
$q1.Cells.Item(15, 1).Value = 13
Set-TextCell $q1 15 2 '014014'
Set-TextCell $q1 15 3 '招商臻选平衡混合A'
Set-TextCell $q1 15 4 '0.84'
Set-TextCell $q1 15 5 '40.88'
Set-TextCell $q1 15 6 '1.00'
Set-TextCell $q1 15 7 '0.0084'
$q1.Cells.Item(15, 8).Value = 9

# 2022-Q1 row 16 (index 14)
$q1.Cells.Item(16, 1).Value = 14
Set-TextCell $q1 16 2 '014015'
Set-TextCell $q1 16 3 '招商臻选平衡混合C'
Set-TextCell $q1 16 4 '0.30'
Set-TextCell $q1 16 5 '40.88'
Set-TextCell $q1 16 6 '1.00'
Set-TextCell $q1 16 7 '0.0030'
$q1.Cells.Item(16, 8).Value = 9

# 2022-Q1 row 17 (index 15)
$q1.Cells.Item(17, 1).Value = 15
Set-TextCell $q1 17 2 '008838'
Set-TextCell $q1 17 3 '德邦量化对冲策略灵活配置混合A'
Set-TextCell $q1 17 4 '0.20'
Set-TextCell $q1 17 5 '38.81'
Set-TextCell $q1 17 6 '1.28'
Set-TextCell $q1 17 7 '0.0026'
$q1.Cells.Item(17, 8).Value = 4

# 2022-Q1 row 18 (index 16)
$q1.Cells.Item(18, 1).Value = 16
Set-TextCell $q1 18 2 '001849'
Set-TextCell $q1 18 3 '前海开源强势共识100强等权重股票'
Set-TextCell $q1 18 4 '0.12'
Set-TextCell $q1 18 5 '92.23'
Set-TextCell $q1 18 6 '1.08'
Set-TextCell $q1 18 7 '0.0013'
$q1.Cells.Item(18, 8).Value = 6

# 2022-Q1 row 19 (index 17)
$q1.Cells.Item(19, 1).Value = 17
Set-TextCell $q1 19 2 '008839'
Set-TextCell $q1 19 3 '德邦量化对冲策略灵活配置混合C'
Set-TextCell $q1 19 4 '0.02'
Set-TextCell $q1 19 5 '38.81'
Set-TextCell $q1 19 6 '1.28'
Set-TextCell $q1 19 7 '0.0003'
$q1.Cells.Item(19, 8).Value = 4

# 2022-Q1 row 20 (index 18)
$q1.Cells.Item(20, 1).Value = 18
Set-TextCell $q1 20 2 '001602'
Set-TextCell $q1 20 3 '鑫元鑫新收益灵活配置混合C'
Set-TextCell $q1 20 4 '0.00'
Set-TextCell $q1 20 5 '86.90'
Set-TextCell $q1 20 6 '4.39'
$q1.Cells.Item(20, 7).Value = 0
$q1.Cells.Item(20, 8).Value = 7

# Re-apply the index-column style (cellXfs 2) to A2:A20, matching the
# sibling fund sheets where the row-number column carries that style.
$q4.Range("A2").Copy()
$q1.Range("A2:A20").PasteSpecial(-4122)  # xlPasteFormats

# ---------------------------------------------------------------------
# 2. Append a brand-new "总计" sheet (gets the next sheetId, 7) after
#    "2022-Q1", holding the totals table with 2022-Q1 prepended.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $q1)
$total.Name = "总计"

$totalHeaders = @("日期", "持有数量(只)", "持有市值(亿元)")
for ($i = 0; $i -lt $totalHeaders.Count; $i++) {
    $total.Cells.Item(1, $i + 2).Value = $totalHeaders[$i]
}
$q4.Range("B1").Copy()
$total.Range("B1:D1").PasteSpecial(-4122)  # xlPasteFormats

# 总计 row 2 (index 0)
$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = '2022-Q1'
$total.Cells.Item(2, 3).Value = 19
$total.Cells.Item(2, 4).Value = 2.2

# 总计 row 3 (index 1)
$total.Cells.Item(3, 1).Value = 1
$total.Cells.Item(3, 2).Value = '2021-Q4'
$total.Cells.Item(3, 3).Value = 4
$total.Cells.Item(3, 4).Value = 2.14

# 总计 row 4 (index 2)
$total.Cells.Item(4, 1).Value = 2
$total.Cells.Item(4, 2).Value = '2021-Q3'
$total.Cells.Item(4, 3).Value = 23
$total.Cells.Item(4, 4).Value = 19.16

# 总计 row 5 (index 3)
$total.Cells.Item(5, 1).Value = 3
$total.Cells.Item(5, 2).Value = '2021-Q2'
$total.Cells.Item(5, 3).Value = 18
$total.Cells.Item(5, 4).Value = 23.42

# 总计 row 6 (index 4)
$total.Cells.Item(6, 1).Value = 4
$total.Cells.Item(6, 2).Value = '2021-Q1'
$total.Cells.Item(6, 3).Value = 17
$total.Cells.Item(6, 4).Value = 5.74

# 总计 row 7 (index 5)
$total.Cells.Item(7, 1).Value = 5
$total.Cells.Item(7, 2).Value = '2020-Q4'
$total.Cells.Item(7, 3).Value = 5
$total.Cells.Item(7, 4).Value = 0.55

$q4.Range("A2").Copy()
$total.Range("A2:A7").PasteSpecial(-4122)  # xlPasteFormats

# Restore the originally-active sheet — adding sheets along the way leaves
# the newest one active/selected, but the source diff doesn't touch bookViews.
$wb.Worksheets.Item("2020-Q4").Activate()
